$d = $word.ActiveDocument

# Disable autocorrect-style quote smartening just in case some code path uses it.
try { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}
try { $word.Options.AutoFormatReplaceQuotes = $false } catch {}

function Replace-Exact {
    param(
        [string]$OldText,
        [string]$NewText,
        [bool]$MatchWholeWord = $false
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0, $MatchWholeWord)
    if (-not $found) {
        throw "Could not find text: $OldText"
    }
    $rng.Text = $NewText
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Exact "Unraveling the Mysteries of the Cosmos" "Unveiling the Enigmatic World of Chemistry"

# ---------------------------------------------------------------------------
# Author line: "Aria Dimitrov" -> "Dr. Caroline Stevens"
# ---------------------------------------------------------------------------
Replace-Exact "Aria Dimitrov" "Dr. Caroline Stevens"

# ---------------------------------------------------------------------------
# Email line: "cosmos_quest@spaceworld" + "." + "net" -> "caroline" + "." + "stevens@schoolmail" + "." + "com"
# ---------------------------------------------------------------------------
Replace-Exact "cosmos_quest@spaceworld" "caroline"
Replace-Exact "net" "stevens@schoolmail.com" $true

# ---------------------------------------------------------------------------
# Body paragraph 1 (font size 24)
# ---------------------------------------------------------------------------
Replace-Exact "The cosmos, an enigmatic tapestry of celestial wonders, has captivated humanity since the dawn of time" "Chemistry, the study of matter, offers a path to unravel the hidden mysteries of the physical world"

Replace-Exact " From ancient astronomers gazing upon the star-studded night sky to modern scientists exploring the outer reaches of the universe, our quest to unravel its mysteries has been an enduring and awe-inspiring pursuit" " It's a science that uncovers the intricacies of tiny particles, revealing an intriguing ballet of atoms and molecules"

Replace-Exact " The cosmos, vast and mysterious, serves as an infinite canvas upon which the grand spectacle of cosmic events unfolds, beckoning us to ponder upon our place in this intricate cosmic symphony" " With each experiment, we peek behind the curtain of the everyday, witnessing the magic of chemical reactions and learning the secrets of substances all around us. Our lives are woven with countless feats of chemistry -- from the food we eat to the medicines that heal us"

Replace-Exact "As we traverse this celestial odyssey, we encounter cosmic phenomena that defy our understanding" "Chemistry allows us to explore the diverse tapestry of materials, from the stardust of diamond to the elasticity of rubber"

Replace-Exact " Supermassive black holes, enigmatic entities lurking at the heart of galaxies, possess a gravitational pull so intense that not even light can escape their clutches" " It's a journey through substances and their interactions, a chronicle of transformations and creations"

Replace-Exact " Neurons, the intricate building blocks of human consciousness, orchestrate a symphony of electrical impulses, enabling us to perceive and navigate the world around us" " Each element tells a story, and each reaction is a chapter in the epic saga of chemistry's symphony"

Replace-Exact " Quantum mechanics, the perplexing realm of the infinitely small, presents us with paradoxes that challenge our very perception of reality" " To understand chemistry is to decipher the language of the physical world, unveiling the dance of particles that defines our existence"

Replace-Exact "These cosmic mysteries, both grand and intricate, ignite a burning curiosity within us" "Chemistry's enigma lies in its universality"

Replace-Exact " They prompt us to delve deeper, to seek answers to questions that have perplexed humanity for millennia" " It transcends species and generations, binding all life together in a web of chemical connections"

Replace-Exact " Our relentless pursuit of understanding has led to profound discoveries, expanding our knowledge of the universe and our place within it" " From the vibrant hues of nature to the intricate machinery of cells, chemistry is a shared language that manifests in countless forms"

Replace-Exact " Yet, as we unravel one enigma, another emerges, beckoning us to continue our exploration, forever captivated by the boundless wonders of the cosmos" " Its universality challenges us to seek patterns in diversity, revealing the interconnectedness of all things. In the realm of chemistry, the ordinary becomes extraordinary as we witness the mundane transformed into a universe of atoms and molecules"

# ---------------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------------
Replace-Exact "Our journey into the cosmos, driven by an insatiable curiosity, has unveiled cosmic wonders that defy comprehension" "Chemistry, the study of matter, unveils the captivating world of substances and reactions"

Replace-Exact " From black holes devouring matter and energy to neurons facilitating the marvels of consciousness, the universe presents us with perplexing phenomena that challenge our understanding" " It takes us on a journey through the diverse tapestry of materials, allowing us to decipher the language of the physical world"

Replace-Exact " Despite the challenges, our unwavering pursuit of knowledge has led to groundbreaking discoveries, propelling us forward in our quest to unravel the intricate mysteries of the cosmos" " Through chemistry, we explore the enigmatic interactions of atoms and molecules, witnessing the epic saga of transformations and creations"

# This span crosses the <w:lastRenderedPageBreak/> run boundary; replacing it
# wholesale both rewrites the text and drops the page-break marker, matching
# the target edit.
Replace-Exact " Our exploration will continue, forever fueled by the allure of the unknown, as we strive to comprehend the enigmatic tapestry of the universe that envelops us" " Its universality binds all life together, inviting us to uncover the interconnectedness of all things. With every experiment, chemistry invites us to question, explore, and discover the hidden marvels of our physical world"

# ---------------------------------------------------------------------------
# Add a trailing empty paragraph after the Summary paragraph (before sectPr).
# ---------------------------------------------------------------------------
$d.Paragraphs.Last.Range.InsertParagraphAfter()

Write-Output "edit complete"
